# Add a new bug-report row (row 17) describing the powerup layering issue.
# Row 17 previously held only empty, pre-formatted cells (styles carried
# over from the table's blank filler rows); it now gets real data that
# matches the look of the other populated rows, so we first clone the
# formatting from row 16 (the previous data row) and then fill in values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone formatting (number formats / fonts / wrap) from the row above so
# the new row visually matches the rest of the table.
$ws.Range("A16:I16").Copy()
$ws.Range("A17:I17").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A17").Value = 45391
$ws.Range("B17").Value = "Annie"
$ws.Range("C17").Value = "P5 - Minor or Feature that would be nice to have"
$ws.Range("D17").Value = "Won't Fix"
$ws.Range("E17").Value = "Annie"
$ws.Range("F17").Value = "Powerups are able to render on top of each other."
$ws.Range("G17").Value = "Expected that there be some jitter with the powerups so that the powerups generate evenly throughout the game board."
$ws.Range("H17").Value = "Powerups are able to layer on top of each other."
$ws.Range("I17").Value = "Play the game for a long enough time, and observe how the powerups generate."
